$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.190.73"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.630.81"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.637.10"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "27.164.95"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.28%  "
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "1.312.59"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "1.767.92"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.801"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +19.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
